$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.131.15"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "2.157.89"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'252.98"
$ws.Range("E5").Value = "  +6.45%  "

$ws.Range("D6").Value = "'0.609"
$ws.Range("E6").Value = "  +1.62%  "

$ws.Range("D7").Value = "'72.93"
$ws.Range("E7").Value = "  +1.62%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.579"
$ws.Range("E9").Value = "  +1.26%  "

$ws.Range("D10").Value = "'39.58"
$ws.Range("E10").Value = "  -0.11%  "

$ws.Range("D11").Value = "'0.0903"
$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("E12").Value = "  +0.68%  "

$ws.Range("D13").Value = "'6.71"
$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("D14").Value = "2.481.53"
$ws.Range("E14").Value = "  -0.25%  "

$ws.Range("D15").Value = "'14.12"
$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("D16").Value = "2.156.34"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").Value = "'0.762"
$ws.Range("E17").Value = "  -1.77%  "

$ws.Range("D18").Value = "42.015.11"
$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("D19").Value = "'0.0000102"
$ws.Range("E19").Value = "  -0.83%  "

$ws.Range("D20").Value = "'70.35"
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("D21").Value = "'5.82"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").Value = "'225.61"
$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("E23").Value = "  -4.79%  "

$ws.Range("D24").Value = "'2.14"
$ws.Range("E24").Value = "  +6.48%  "

$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").Value = "'10.40"
$ws.Range("E26").Value = "  -2.74%  "

$ws.Range("D27").Value = "'3.31"
$ws.Range("E27").Value = "  +1.16%  "

$ws.Range("E28").Value = "  +6.07%  "

$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").Value = "'36.61"
$ws.Range("E30").Value = "  +11.12%  "

$ws.Range("D31").Value = "'167.95"
$ws.Range("E31").Value = "  -1.62%  "

$ws.Range("D32").Value = "'19.88"
$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("D33").Value = "'0.0802"
$ws.Range("E33").Value = "  +4.30%  "

$ws.Range("D34").Value = "'5.09"
$ws.Range("E34").Value = "  -4.26%  "

$ws.Range("D35").Value = "'0.120"
$ws.Range("E35").Value = "  -0.52%  "

$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = "  +4.84%  "

$ws.Range("D37").Value = "'4.22"
$ws.Range("E37").Value = "  -1.38%  "

$ws.Range("D38").Value = "'0.0328"
$ws.Range("E38").Value = "  +8.26%  "

$ws.Range("D39").Value = "'11.75"
$ws.Range("E39").Value = "  -2.22%  "

$ws.Range("D40").Value = "'2.04"
$ws.Range("E40").Value = "  -2.47%  "

$ws.Range("E41").Value = "  +3.59%  "

$ws.Range("D42").Value = "'58.39"
$ws.Range("E42").Value = "  -0.82%  "

$ws.Range("D43").Value = "'5.10"
$ws.Range("E43").Value = "  -4.66%  "

$ws.Range("D44").Value = "'101.73"
$ws.Range("E44").Value = "  +5.01%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'8.19"
$ws.Range("E45").Value = "  -3.32%  "

$ws.Range("B46").Value = "WOONetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D46").Value = "'0.458"
$ws.Range("E46").Value = "  +15.01%  "

$ws.Range("D47").Value = "'0.0959"
$ws.Range("E47").Value = "  -0.57%  "

$ws.Range("D48").Value = "'2.38"
$ws.Range("E48").Value = "  +9.41%  "

$ws.Range("D49").Value = "'1.08"
$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D50").Value = "'1.12"
$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("E51").Value = "  +0.58%  "

